$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: update existing rows 189-191 values
$ws.Cells.Item(189,2).Value = 7
$ws.Cells.Item(190,2).Value = 10
$ws.Cells.Item(191,2).Value = 11

# Step 2: bulk-insert 39 rows at position 195, pushing the old SUM row (195) down to 234
$ws.Rows("195:233").Insert(-4121, -4163)

# Step 3: row 192 (pre-existing blank slot) - style A4 B6 C4 D4 (matches rows 178-191 pattern)
$ws.Range("A191:D191").Copy()
$ws.Range("A192:D192").PasteSpecial(-4122)
$ws.Cells.Item(192,1).Value = 44976
$ws.Cells.Item(192,2).Value = 5

# Step 4: rows 194-200 - style A4 B5 C4 D4 (A/C/D from row176, B from row9)
$ws.Range("A176").Copy()
$ws.Range("A194:A200").PasteSpecial(-4122)
$ws.Range("C176:D176").Copy()
$ws.Range("C194:D200").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B194:B200").PasteSpecial(-4122)
$ws.Cells.Item(194,1).Value = 44977
$ws.Cells.Item(194,2).Value = 1
$ws.Cells.Item(195,1).Value = 44978
$ws.Cells.Item(195,2).Value = 5
$ws.Cells.Item(196,1).Value = 44979
$ws.Cells.Item(196,2).Value = 2
$ws.Cells.Item(197,1).Value = 44980
$ws.Cells.Item(197,2).Value = 5
$ws.Cells.Item(198,1).Value = 44981
$ws.Cells.Item(198,2).Value = 3
$ws.Cells.Item(199,1).Value = 44982
$ws.Cells.Item(199,2).Value = 10
$ws.Cells.Item(200,1).Value = 44983
$ws.Cells.Item(200,2).Value = 5

# Step 5: rows 202-232 - style A4 B5 C5 D5 (matches row9-175 dominant pattern)
$ws.Range("A9:D9").Copy()
$ws.Range("A202:D208").PasteSpecial(-4122)
$ws.Range("A210:D216").PasteSpecial(-4122)
$ws.Range("A218:D224").PasteSpecial(-4122)
$ws.Range("A226:D232").PasteSpecial(-4122)
$ws.Cells.Item(202,1).Value = 44984
$ws.Cells.Item(202,2).Value = 3
$ws.Cells.Item(203,1).Value = 44985
$ws.Cells.Item(203,2).Value = 2
$ws.Cells.Item(204,1).Value = 44986
$ws.Cells.Item(204,2).Value = 2
$ws.Cells.Item(205,1).Value = 44987
$ws.Cells.Item(205,2).Value = 3
$ws.Cells.Item(206,1).Value = 44988
$ws.Cells.Item(206,2).Value = 3
$ws.Cells.Item(207,1).Value = 44989
$ws.Cells.Item(207,2).Value = 9
$ws.Cells.Item(208,1).Value = 44990
$ws.Cells.Item(208,2).Value = 7

# Dates for empty rows 210-216, 218-224, 226-232
$ws.Cells.Item(210,1).Value = 44991
$ws.Cells.Item(211,1).Value = 44992
$ws.Cells.Item(212,1).Value = 44993
$ws.Cells.Item(213,1).Value = 44994
$ws.Cells.Item(214,1).Value = 44995
$ws.Cells.Item(215,1).Value = 44996
$ws.Cells.Item(216,1).Value = 44997
$ws.Cells.Item(218,1).Value = 44998
$ws.Cells.Item(219,1).Value = 44999
$ws.Cells.Item(220,1).Value = 45000
$ws.Cells.Item(221,1).Value = 45001
$ws.Cells.Item(222,1).Value = 45002
$ws.Cells.Item(223,1).Value = 45003
$ws.Cells.Item(224,1).Value = 45004
$ws.Cells.Item(226,1).Value = 45005
$ws.Cells.Item(227,1).Value = 45006
$ws.Cells.Item(228,1).Value = 45007
$ws.Cells.Item(229,1).Value = 45008
$ws.Cells.Item(230,1).Value = 45009
$ws.Cells.Item(231,1).Value = 45010
$ws.Cells.Item(232,1).Value = 45011

# Step 6: update SUM row (now at 234)
$ws.Cells.Item(234,2).Formula = "=SUM(B2:B232)"
$ws.Cells.Item(234,3).Formula = "=B234/23"
$ws.Cells.Item(234,4).Formula = "=23*7"
